$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as text (matching the workbook convention)
# instead of auto-converting them to numbers.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D11",
    "D12",
    "D13",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D39",
    "D44",
    "D46"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "70.649.12"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.560.55"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "587.44"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "189.51"
$ws.Range("E6").Value = "  +4.31%  "
$ws.Range("D7").Value = "3.552.55"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "0.215"
$ws.Range("E10").Value = "  +10.38%  "
$ws.Range("D11").Value = "0.644"
$ws.Range("D12").Value = "54.10"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "0.0000310"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "4.123.99"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "70.659.81"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "3.586.87"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "12.68"
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("D19").Value = "18.91"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "564.18"
$ws.Range("E20").Value = "  +5.01%  "
$ws.Range("B21").Value = "TRON"
$ws.Range("C21").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D21").Value = "0.120"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "17.89"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "4.87"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "93.77"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").Value = "32.30"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "7.06"
$ws.Range("D32").Value = "12.16"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "3.98"
$ws.Range("E33").Value = "  +30.05%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "63.16"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "530.21"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "38.06"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("D40").Value = "3.665.78"
$ws.Range("E40").Value = "  +10.62%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "0.0₃0786"
$ws.Range("E42").Value = "  +4.16%  "
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("D44").Value = "0.137"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("D46").Value = "3.48"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +8.33%  "
